$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so that values such as
# "1.191" or "27.253.75" are stored as literal text, not parsed as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.253.75"
$ws.Range("E2").Value = "  -0.76%  "

# Row 3
$ws.Range("D3").Value = "1.783.86"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.67%  "

# Row 5
$ws.Range("D5").Value = "335.67"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +1.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.35%  "

# Row 8
$ws.Range("D8").Value = "0.3411"
$ws.Range("E8").Value = "  -2.02%  "

# Row 9
$ws.Range("D9").Value = "48.24"
$ws.Range("E9").Value = "  -2.37%  "

# Row 10
$ws.Range("D10").Value = "1.191"
$ws.Range("E10").Value = "  -2.55%  "

# Row 11
$ws.Range("D11").Value = "0.07426"
$ws.Range("E11").Value = "  -3.47%  "

# Row 12
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("D13").Value = "21.57"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("D14").Value = "6.397"
$ws.Range("E14").Value = "  -2.83%  "

# Row 15
$ws.Range("D15").Value = "1.782.19"
$ws.Range("E15").Value = "  +0.28%  "

# Row 16
$ws.Range("D16").Value = "7.024"
$ws.Range("E16").Value = "  -2.74%  "

# Row 17
$ws.Range("E17").Value = "  -2.24%  "

# Row 18
$ws.Range("D18").Value = "0.06675"
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
$ws.Range("D19").Value = "84.09"
$ws.Range("E19").Value = "  -0.99%  "

# Row 20
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.94%  "

# Row 21
$ws.Range("D21").Value = "6.525"
$ws.Range("E21").Value = "  +2.19%  "

# Row 22
$ws.Range("D22").Value = "17.23"
$ws.Range("E22").Value = "  -1.61%  "

# Row 23
$ws.Range("D23").Value = "27.194.72"
$ws.Range("E23").Value = "  -0.94%  "

# Row 24
$ws.Range("D24").Value = "12.39"
$ws.Range("E24").Value = "  -5.85%  "

# Row 25
$ws.Range("D25").Value = "2.415"
$ws.Range("E25").Value = "  -2.46%  "

# Row 26
$ws.Range("D26").Value = "1.501"
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("D27").Value = "2.519"
$ws.Range("E27").Value = "  -0.40%  "

# Row 28
$ws.Range("D28").Value = "21.09"
$ws.Range("E28").Value = "  +4.44%  "

# Row 29
$ws.Range("D29").Value = "152.47"
$ws.Range("E29").Value = "  -0.45%  "

# Row 30
$ws.Range("D30").Value = "1.982.60"
$ws.Range("E30").Value = "  +0.48%  "

# Row 31
$ws.Range("D31").Value = "132.54"
$ws.Range("E31").Value = "  -1.99%  "

# Row 32
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  -0.97%  "

# Row 33
$ws.Range("D33").Value = "5.992"
$ws.Range("E33").Value = "  -4.65%  "

# Row 34
$ws.Range("D34").Value = "0.08582"
$ws.Range("E34").Value = "  -1.89%  "

# Row 35
$ws.Range("D35").Value = "13.01"
$ws.Range("E35").Value = "  -3.05%  "

# Row 36
$ws.Range("D36").Value = "1.650"
$ws.Range("E36").Value = "  -4.19%  "

# Row 37
$ws.Range("D37").Value = "5.394"
$ws.Range("E37").Value = "  -4.14%  "

# Row 38
$ws.Range("D38").Value = "0.6802"
$ws.Range("E38").Value = "  +0.39%  "

# Row 39
$ws.Range("D39").Value = "0.06319"
$ws.Range("E39").Value = "  -2.14%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02337"
$ws.Range("E40").Value = "  -2.88%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "8.765"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.2175"
$ws.Range("E42").Value = "  -3.23%  "

# Row 43
$ws.Range("D43").Value = "1.245"
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
$ws.Range("D44").Value = "14.38"
$ws.Range("E44").Value = "  -1.88%  "

# Row 45
$ws.Range("D45").Value = "1.004"
$ws.Range("E45").Value = "  +0.92%  "

# Row 46
$ws.Range("D46").Value = "0.6343"
$ws.Range("E46").Value = "  -1.74%  "

# Row 47
$ws.Range("D47").Value = "3.842"
$ws.Range("E47").Value = "  -3.25%  "

# Row 48
$ws.Range("D48").Value = "2.113"
$ws.Range("E48").Value = "  -1.64%  "

# Row 49
$ws.Range("D49").Value = "128.58"
$ws.Range("E49").Value = "  -2.10%  "

# Row 50
$ws.Range("D50").Value = "0.07162"
$ws.Range("E50").Value = "  -2.58%  "

# Row 51
$ws.Range("D51").Value = "78.93"
$ws.Range("E51").Value = "  -1.61%  "
